# Scheduled market-price refresh for Maduin_Profits (crafting leve profit sheets).
# Updates cached currentAveragePrice*/LevePrice*/LeveProfit* columns (H:N) per row
# to the latest Universalis snapshot values pulled by the runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12 (Leve Item ID 5515)
$ws.Range("H12").Value = 1025.2
$ws.Range("I12").Value = 1326.1818
$ws.Range("J12").Value = 197.5
$ws.Range("K12").Value = 1326.1818
$ws.Range("L12").Value = 197.5
$ws.Range("M12").Value = -1156.1818
$ws.Range("N12").Value = -537.5

# Row 13 (Leve Item ID 2144)
$ws.Range("H13").Value = 1000
$ws.Range("J13").Value = 1000
$ws.Range("L13").Value = 1000
$ws.Range("N13").Value = -1338

# Row 74 (Leve Item ID 5507)
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = ""

# Row 77 (Leve Item ID 5507)
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = ""

# Row 106 (Leve Item ID 19903)
$ws.Range("H106").Value = 400
$ws.Range("I106").Value = 400
$ws.Range("K106").Value = 400
$ws.Range("M106").Value = 231

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 1054.8334
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = ""

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 2950
$ws.Range("I137").Value = 2950
$ws.Range("K137").Value = 8850
$ws.Range("M137").Value = -6300

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 5765.952
$ws.Range("I138").Value = 4814.8
$ws.Range("J138").Value = 6063.1875
$ws.Range("K138").Value = 14444.4
$ws.Range("L138").Value = 18189.5625
$ws.Range("M138").Value = -9304.400000000001
$ws.Range("N138").Value = -28469.5625

$ws = $wb.Worksheets.Item("ARM")
# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 2874.5
$ws.Range("I110").Value = 2999.4285
$ws.Range("K110").Value = 2999.4285
$ws.Range("M110").Value = -954.4285

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 740
$ws.Range("I122").Value = 716.6667
$ws.Range("J122").Value = 775
$ws.Range("K122").Value = 2150.0001
$ws.Range("L122").Value = 2325
$ws.Range("M122").Value = 299.9998999999998
$ws.Range("N122").Value = -7225

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 218.33333
$ws.Range("I22").Value = 77.5
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 77.5
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 95.5
$ws.Range("N22").Value = -846

# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 2649.7
$ws.Range("I86").Value = 2716.2666
$ws.Range("J86").Value = 2450
$ws.Range("K86").Value = 2716.2666
$ws.Range("L86").Value = 2450
$ws.Range("M86").Value = -1593.2666
$ws.Range("N86").Value = -4696

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 2649.7
$ws.Range("I89").Value = 2716.2666
$ws.Range("J89").Value = 2450
$ws.Range("K89").Value = 13581.333
$ws.Range("L89").Value = 12250
$ws.Range("M89").Value = -7965.332999999999
$ws.Range("N89").Value = -23482

# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 3510.5
$ws.Range("I107").Value = 1350
$ws.Range("J107").Value = 5671
$ws.Range("K107").Value = 1350
$ws.Range("L107").Value = 5671
$ws.Range("M107").Value = 570
$ws.Range("N107").Value = -9511

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 951.1177
$ws.Range("I134").Value = 891.3333
$ws.Range("K134").Value = 2673.9999
$ws.Range("M134").Value = -138.9998999999998

$ws = $wb.Worksheets.Item("CRP")
# Row 4 (Leve Item ID 3742)
$ws.Range("H4").Value = 15983.6875
$ws.Range("J4").Value = 18710.77
$ws.Range("L4").Value = 18710.77
$ws.Range("N4").Value = -18934.77

# Row 14 (Leve Item ID 1998)
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = ""

# Row 96 (Leve Item ID 18193)
$ws.Range("H96").Value = 19108.166
$ws.Range("J96").Value = 19108.166
$ws.Range("L96").Value = 19108.166
$ws.Range("N96").Value = -24600.166

$ws = $wb.Worksheets.Item("CUL")
# Row 17 (Leve Item ID 4640)
$ws.Range("H17").Value = 7625.25
$ws.Range("I17").Value = 10000
$ws.Range("K17").Value = 30000
$ws.Range("M17").Value = -29831

# Row 26 (Leve Item ID 4746)
$ws.Range("H26").Value = 70
$ws.Range("I26").Value = 69.333336
$ws.Range("J26").Value = 71
$ws.Range("K26").Value = 208.000008
$ws.Range("L26").Value = 213
$ws.Range("M26").Value = 79.99999199999999
$ws.Range("N26").Value = -789

# Row 33 (Leve Item ID 4867)
$ws.Range("H33").Value = 116
$ws.Range("I33").Value = 144.5
$ws.Range("J33").Value = 2
$ws.Range("K33").Value = 867
$ws.Range("L33").Value = 12
$ws.Range("M33").Value = -584
$ws.Range("N33").Value = -578

# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 3215.6428
$ws.Range("I68").Value = 3878.25
$ws.Range("J68").Value = 2332.1667
$ws.Range("K68").Value = 11634.75
$ws.Range("L68").Value = 6996.500100000001
$ws.Range("M68").Value = -10823.75
$ws.Range("N68").Value = -8618.500100000001

# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 3215.6428
$ws.Range("I71").Value = 3878.25
$ws.Range("J71").Value = 2332.1667
$ws.Range("K71").Value = 34904.25
$ws.Range("L71").Value = 20989.5003
$ws.Range("M71").Value = -30848.25
$ws.Range("N71").Value = -29101.5003

# Row 87 (Leve Item ID 12864)
$ws.Range("H87").Value = 250
$ws.Range("I87").Value = 250
$ws.Range("K87").Value = 750
$ws.Range("M87").Value = 498

# Row 90 (Leve Item ID 12864)
$ws.Range("H90").Value = 250
$ws.Range("I90").Value = 250
$ws.Range("K90").Value = 2250
$ws.Range("M90").Value = 3990

# Row 92 (Leve Item ID 19841)
$ws.Range("H92").Value = 642.8570999999999

# Row 121 (Leve Item ID 27878)
$ws.Range("H121").Value = 1947.5
$ws.Range("I121").Value = 395
$ws.Range("K121").Value = 1185
$ws.Range("M121").Value = 125

# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 990
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").Value = ""

# Row 137 (Leve Item ID 44088)
$ws.Range("H137").Value = 3850.5715
$ws.Range("I137").Value = 3749.75
$ws.Range("K137").Value = 11249.25
$ws.Range("M137").Value = -6149.25

$ws = $wb.Worksheets.Item("LTW")
# Row 2 (Leve Item ID 2631)
$ws.Range("H2").Value = 35000
$ws.Range("J2").Value = 50000
$ws.Range("L2").Value = 50000
$ws.Range("N2").Value = -50224

# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 2258.3635
$ws.Range("I22").Value = 1718.4
$ws.Range("J22").Value = 2708.3333
$ws.Range("K22").Value = 1718.4
$ws.Range("L22").Value = 2708.3333
$ws.Range("M22").Value = -1423.4
$ws.Range("N22").Value = -3298.3333

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 2258.3635
$ws.Range("I27").Value = 1718.4
$ws.Range("J27").Value = 2708.3333
$ws.Range("K27").Value = 1718.4
$ws.Range("L27").Value = 2708.3333
$ws.Range("M27").Value = -1611.4
$ws.Range("N27").Value = -2922.3333

# Row 43 (Leve Item ID 4314)
$ws.Range("H43").Value = 26500
$ws.Range("J43").Value = 3000
$ws.Range("L43").Value = 3000
$ws.Range("N43").Value = -3386

$ws = $wb.Worksheets.Item("WVR")
# Row 2 (Leve Item ID 3307)
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = ""

# Row 24 (Leve Item ID 3561)
$ws.Range("H24").Value = 5000
$ws.Range("J24").Value = 5000
$ws.Range("L24").Value = 5000
$ws.Range("N24").Value = -5460

# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 1221.1538
$ws.Range("I107").Value = 946
$ws.Range("K107").Value = 2838
$ws.Range("M107").Value = -918

# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 2394.75
$ws.Range("I122").Value = 2394.75
$ws.Range("K122").Value = 7184.25
$ws.Range("M122").Value = -4734.25

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 1245.625
$ws.Range("I136").Value = 852.1429000000001
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 2556.4287
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -6.42870000000039
$ws.Range("N136").Value = -17100
